$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: replace the old header row with a single title cell ---
# Clear the whole header row first (old headers spanned A1:L1)
$ws.Range("A1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# --- Rows 2-10: columns get reshuffled ---
# Old layout: A=Day B=Time C=ModuleCode D=ModuleTitle E=Hours F=ClassType
#             G=Lecturer H=Room I=Block J=Group K=Level L=Course
# New layout: A=Day B=Time C=Hours D=ModuleCode E=ModuleTitle F=ClassType
#             G=Lecturer H=Group I=Block J=Room
# (Level and Course columns are removed entirely)

$data = @(
    @("SUN", "7:00-9:30",   2.5, "5CS024", "Collaborative Development",                 "Workshop", "Mr. Raj Shrestha",    "L5CG14",          "WLV", "TR-01 Dudley"),
    @("SUN", "10:00-12:00", 2,   "5CS020", "Distributed and Cloud Systems Programming", "Tutorial", "Mr. Sumanta Silwal",  "L5CG14",          "HCK", "TR-06 Nagpokhari"),
    @("MON", "9:30-12:00",  2.5, "5CS020", "Distributed and Cloud Systems Programming", "Workshop", "Mr. Sumanta Silwal",  "L5CG14",          "HCK", "Lab-03 Gahanapokhari"),
    @("MON", "13:00-15:30", 2.5, "5CS022", "Human Computer Interaction",                "Workshop", "Mr. Dipesh Shrestha", "L5CG14",          "WLV", "SR-02 Bilston"),
    @("TUE", "7:00-9:00",   2,   "5CS024", "Collaborative Development",                 "Lecture",  "Mr. Raj Shrestha",    "L5CG(12+13+14)",  "WLV", "LT-03 Walsall"),
    @("WED", "7:00-9:00",   2,   "5CS022", "Human Computer Interaction",                "Lecture",  "Mr. Ayush Shakya",    "L5CG(12+13+14)",  "WLV", "LT-01 Wulfruna"),
    @("THU", "7:00-9:00",   2,   "5CS024", "Collaborative Development",                 "Tutorial", "Mr. Raj Shrestha",    "L5CG14",          "WLV", "TR-03 Westbromwich"),
    @("THU", "9:30-11:30",  2,   "5CS020", "Distributed and Cloud Systems Programming", "Lecture",  "Mr. Sumanta Silwal",  "L5CG(12+13+14)",  "WLV", "LT-01 Wulfruna"),
    @("FRI", "9:00-11:00",  2,   "5CS022", "Human Computer Interaction",                "Tutorial", "Mr. Dipesh Shrestha", "L5CG14",          "WLV", "SR-01 Bantok")
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

# --- Remove the now-unused Level/Course columns (old K and L) ---
$ws.Range("K1:L10").Delete()
